$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 300
$ws1.Range("F3").Value = 1201
$ws1.Range("F4").Value = 16834
$ws1.Range("F10").Value = 220
$ws1.Range("F11").Value = 128
$ws1.Range("F12").Value = 11669
$ws1.Range("F14").Value = 1342
$ws1.Range("F15").Value = 4629
$ws1.Range("F16").Value = 454
$ws1.Range("F17").Value = 404
$ws1.Range("F18").Value = 69
$ws1.Range("F19").Value = 895

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 300
$ws4.Range("F4").Value = 1201
$ws4.Range("F5").Value = 16834
$ws4.Range("F11").Value = 220
$ws4.Range("F12").Value = 128
$ws4.Range("F15").Value = 11669
$ws4.Range("F17").Value = 1342
$ws4.Range("F18").Value = 4629
$ws4.Range("F19").Value = 454
$ws4.Range("F20").Value = 404
$ws4.Range("F21").Value = 69
$ws4.Range("F22").Value = 895
